# Fixed exception naming for RunAllTests.
# Rename the "AppEx" exception category to "SystemException" and align the
# "BRE" category name to "BusinessException" in both the Tests and Result
# sheets' data-validation dropdown lists, widen the validated ranges to the
# full column, fix up the one cell that used the old "AppEx" status, and
# restore the "Tests" tab as the active/selected sheet (swapping places
# with "Result").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tests")
$ws2 = $wb.Worksheets.Item("Result")

# --- Tests sheet: fix the one cell that still says "AppEx" -> "SystemException"
$ws1.Range("B5").Value = "SystemException"

# --- Tests sheet: data validation list + widened range
$ws1.Range("B2:B16").Validation.Delete() | Out-Null
$ws1.Range("B2:B1048576").Validation.Add(3, 1, 1, '"Success,BusinessException,SystemException"') | Out-Null

# --- Result sheet: data validation list + widened range
$ws2.Range("B2:B17").Validation.Delete() | Out-Null
$ws2.Range("B2:B1048576").Validation.Add(3, 1, 1, '"Success,BusinessException,SystemException"') | Out-Null

# --- Selection bookkeeping: Result used to be the active tab with D15
# selected; Tests is now active with F12 selected, and Result keeps a
# stored selection of D14.
$ws2.Range("D14").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F12").Select() | Out-Null
